$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 9) following the existing "label / ${placeholder}" pattern
$ws.Range("F9").Value = "image"
$ws.Range("G9").Value = "`${image}"

# Move the active selection the way Excel would after typing into G9 and pressing Enter
$ws.Range("G10").Select()
